$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: clone header style (s=1) onto new columns, then set text ---
$ws.Range("C1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Indian"
$ws.Range("D1").Value = "Foreign"
$ws.Range("E1").Value = "discharged_y"
$ws.Range("F1").Value = "deaths_y"
$ws.Range("G1").Value = "helpline"

# --- Data rows ---
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "8662410978"
$ws.Range("G2").Style = "Normal"

$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "9436055743"
$ws.Range("G3").Style = "Normal"

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "6913347770"
$ws.Range("G4").Style = "Normal"

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "104"
$ws.Range("G5").Style = "Normal"

$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = "077122-35091"

$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "104"
$ws.Range("G7").Style = "Normal"

$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "104"
$ws.Range("G8").Style = "Normal"

$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 14
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "8558893911"
$ws.Range("G9").Style = "Normal"

$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "104"
$ws.Range("G10").Style = "Normal"

$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "104"
$ws.Range("G11").Style = "Normal"

$ws.Range("C12").Value = 15
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "104"
$ws.Range("G12").Style = "Normal"

$ws.Range("C13").Value = 26
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = "0471-2552056"

$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = "0755-2527177"

$ws.Range("C15").Value = 49
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = "020-26127394"

$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "3852411668"
$ws.Range("G16").Style = "Normal"

$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "108"
$ws.Range("G17").Style = "Normal"

$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "102"
$ws.Range("G18").Style = "Normal"

$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "7005539653"
$ws.Range("G19").Style = "Normal"

$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "9439994859"
$ws.Range("G20").Style = "Normal"

$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "104"
$ws.Range("G21").Style = "Normal"

$ws.Range("C22").Value = 15
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = "0141-2225624"

$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "104"
$ws.Range("G23").Style = "Normal"

$ws.Range("C24").Value = 3
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = "044-29510500"

$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "104"
$ws.Range("G25").Style = "Normal"

$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = "0381-2315879"

$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "104"
$ws.Range("G27").Style = "Normal"

$ws.Range("C28").Value = 22
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 9
$ws.Range("F28").Value = 0
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "18001805145"
$ws.Range("G28").Style = "Normal"

$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "3323412600"
$ws.Range("G29").Style = "Normal"

$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = "03192-232102"

$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "9779558282"
$ws.Range("G31").Style = "Normal"

$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "104"
$ws.Range("G32").Style = "Normal"

$ws.Range("C33").Value = 16
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = 5
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = "011-22307154"

$ws.Range("C34").Value = 4
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = "1912520982 (Jammu), 1942440283 (Kashmir) "

$ws.Range("C35").Value = 10
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "1982256462"
$ws.Range("G35").Style = "Normal"

$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "104"
$ws.Range("G36").Style = "Normal"

$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "104"
$ws.Range("G37").Style = "Normal"
